# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback
# DateTime" (H2) timestamps on the per-language handback status sheets.
# These cells store plain text (not real Excel dates), so we keep them as
# text strings on write.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 06:52:31"
$wsZhCn.Range("H2").Value = "2016-03-20 06:52:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 06:52:34"
$wsDeDe.Range("H2").Value = "2016-03-20 06:52:55"
